# Apply last week's changes:
#  - Remove the section-divider slide that only had an empty "Title 1" placeholder (old slide 2)
#  - Remove the empty chart-overview slide (old slide 4)
#  - Update the date placeholder text on the "title"/date slide from 20/02/2020 to 21/02/2020
#  - Remove the "This is the section" section-header slide (old slide 6)

$p = $ppt.ActivePresentation

# Update the date text first, while slide indices still match the original deck
# (old slide 5 -> "Text Placeholder 3" holds the date string).
$dateSlide = $p.Slides.Item(5)
$dateSlide.Shapes.Item(3).TextFrame.TextRange.Text = "21/02/2020"

# Delete slides from the end towards the start so earlier indices remain valid.
$p.Slides.Item(6).Delete()   # "This is the section" slide
$p.Slides.Item(4).Delete()   # empty chart placeholder slide
$p.Slides.Item(2).Delete()   # empty ctrTitle-only slide
